$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 324.42856
$ws.Range("I33").Value = 391.0909
$ws.Range("K33").Value = 391.0909
$ws.Range("M33").Value = -162.0909

$ws.Range("H43").Value = 4032.6667
$ws.Range("J43").Value = 5378.2
$ws.Range("L43").Value = 5378.2
$ws.Range("N43").Value = -5516.2

$ws.Range("H88").Value = 3090.5
$ws.Range("I88").Value = 2683.8572
$ws.Range("J88").Value = 3659.8
$ws.Range("K88").Value = 2683.8572
$ws.Range("L88").Value = 3659.8
$ws.Range("M88").Value = -2277.8572
$ws.Range("N88").Value = -4471.8

$ws.Range("H91").Value = 3090.5
$ws.Range("I91").Value = 2683.8572
$ws.Range("J91").Value = 3659.8
$ws.Range("K91").Value = 2683.8572
$ws.Range("L91").Value = 3659.8
$ws.Range("M91").Value = -1279.8572
$ws.Range("N91").Value = -6467.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1872.75
$ws.Range("I45").Value = 1872.75
$ws.Range("K45").Value = 1872.75
$ws.Range("M45").Value = -1495.75

$ws.Range("H61").Value = 1639.8948
$ws.Range("I61").Value = 1597.5294
$ws.Range("K61").Value = 1597.5294
$ws.Range("M61").Value = -1385.5294

$ws.Range("H92").Value = 45774.75
$ws.Range("J92").Value = 45774.75
$ws.Range("L92").Value = 45774.75
$ws.Range("N92").Value = -50766.75

$ws.Range("H97").Value = 5300
$ws.Range("I97").Value = 600
$ws.Range("J97").Value = 10000
$ws.Range("K97").Value = 600
$ws.Range("L97").Value = 10000
$ws.Range("M97").Value = -104
$ws.Range("N97").Value = -10992

$ws.Range("H102").Value = 1662.375
$ws.Range("I102").Value = 1757
$ws.Range("K102").Value = 1757
$ws.Range("M102").Value = -135

$ws.Range("H122").Value = 1447.375
$ws.Range("I122").Value = 1458.1305
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 4374.3915
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -1924.3915
$ws.Range("N122").Value = -8500

$ws.Range("H136").Value = 1639.8948
$ws.Range("I136").Value = 1597.5294
$ws.Range("K136").Value = 4792.5882
$ws.Range("M136").Value = -2242.5882

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1624.25
$ws.Range("I31").Value = 1332.3334
$ws.Range("K31").Value = 1332.3334
$ws.Range("M31").Value = -1037.3334

$ws.Range("H34").Value = 1624.25
$ws.Range("I34").Value = 1332.3334
$ws.Range("K34").Value = 1332.3334
$ws.Range("M34").Value = -1130.3334

$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996

$ws.Range("H82").Value = 23300
$ws.Range("I82").Value = 23300
$ws.Range("K82").Value = 23300
$ws.Range("M82").Value = -22939

$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984

$ws.Range("H85").Value = 23300
$ws.Range("I85").Value = 23300
$ws.Range("K85").Value = 23300
$ws.Range("M85").Value = -22052

$ws.Range("H94").Value = 1115.6666
$ws.Range("J94").Value = 673.5
$ws.Range("L94").Value = 673.5
$ws.Range("N94").Value = -1575.5

$ws.Range("H107").Value = 1916
$ws.Range("J107").Value = 2188
$ws.Range("L107").Value = 2188
$ws.Range("N107").Value = -6028

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 16499.875
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 18285.572
$ws.Range("K70").Value = 12000
$ws.Range("L70").Value = 54856.716
$ws.Range("M70").Value = -11685
$ws.Range("N70").Value = -55486.716

$ws.Range("H73").Value = 16499.875
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 18285.572
$ws.Range("K73").Value = 12000
$ws.Range("L73").Value = 54856.716
$ws.Range("M73").Value = -10908
$ws.Range("N73").Value = -57040.716

$ws.Range("H137").Value = 7229.6665
$ws.Range("J137").Value = 8632
$ws.Range("L137").Value = 25896
$ws.Range("N137").Value = -36096

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1794.7587
$ws.Range("I122").Value = 1285.96
$ws.Range("J122").Value = 4974.75
$ws.Range("K122").Value = 3857.88
$ws.Range("L122").Value = 14924.25
$ws.Range("M122").Value = -1407.88
$ws.Range("N122").Value = -19824.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2069.75
$ws.Range("I46").Value = 2093
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 2093
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -1905
$ws.Range("N46").Value = -2376

$ws.Range("H122").Value = 3628.25
$ws.Range("I122").Value = 3171.1667
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 9513.500100000001
$ws.Range("L122").Value = 14998.5
$ws.Range("M122").Value = -7063.500100000001
$ws.Range("N122").Value = -19898.5

$ws.Range("H132").Value = 3080.8
$ws.Range("I132").Value = 3310.7
$ws.Range("J132").Value = 2621
$ws.Range("K132").Value = 9932.099999999999
$ws.Range("L132").Value = 7863
$ws.Range("M132").Value = -7402.099999999999
$ws.Range("N132").Value = -12923

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 11995
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

$ws.Range("H82").Value = 23200
$ws.Range("I82").Value = 23200
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 23200
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -22817
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 23200
$ws.Range("I85").Value = 23200
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 23200
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -21874
$ws.Range("N85").ClearContents()

$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()

$ws.Range("H141").Value = 49000
$ws.Range("J141").Value = 48000
$ws.Range("L141").Value = 48000
$ws.Range("N141").Value = -58360
